# update PANELS materials for wet weight
#
# Glossary sheet: the old "temperature_c" glossary row (row 30) is repurposed
# to hold the "topdown_photos_notes" definition (what used to live in row 31),
# and row 31 is repurposed to introduce the new "wet_weight_g" field.
$wb = $excel.ActiveWorkbook

$glossary = $wb.Worksheets.Item("glossary")
$glossary.Range("A30").Value = "topdown_photos_notes"
$glossary.Range("B30").Value = "Any additional notes regarding observations, context, or concerns about the data."
$glossary.Range("C30").Value = "text"
$glossary.Range("F30").Value = "topdown photos"

$glossary.Range("A31").Value = "wet_weight_g"
$glossary.Range("B31").Value = "The total wet mass of the entire community that was scrapped off of the panel"
$glossary.Range("C31").Value = "numeric"
$glossary.Range("F31").Value = "eDNA"

# site metadata sheet: the temperature_c column is no longer collected there,
# so remove it (columns to its right shift left).
$siteMetadata = $wb.Worksheets.Item("site metadata")
$siteMetadata.Range("K1").EntireColumn.Delete()

# eDNA sheet: add the new wet_weight_g column, ahead of the existing eDNA_notes
# column (which shifts from D to E).
$eDNA = $wb.Worksheets.Item("eDNA")
$eDNA.Range("D1").EntireColumn.Insert()
$eDNA.Range("D1").Value = "wet_weight_g"
$eDNA.Columns.Item(4).ColumnWidth = 11.75
